$d = $word.ActiveDocument

# Merge the "Versi" + "on" runs into a single "Version" run.
$d.Content.Find.Execute("Version", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Version", 2)

# Change " 2" to " 1." (the new period belongs in this run, before the bookmark).
$d.Content.Find.Execute(" 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, " 1.", 2)

# Remove the now-duplicate trailing "." run that sits after the bookmark.
$r = $d.Range(10, 11)
$r.Delete()
